$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 23-43 (the extra data rows that get removed)
$ws.Range("A23:B43").EntireRow.Delete()

# Update column A (rows 2-22) to all be "Alex"
$ws.Range("A2:A22").Value = "Alex"

# Update column B (rows 2-22) to all be "Utrecht"
$ws.Range("B2:B22").Value = "Utrecht"

# Update the selection / view
$ws.Range("B12").Select()
